$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F ("Current balance"), shifting existing
# columns F:T to G:U, to make room for the new "Previous balance" column.
$ws.Columns("F").Insert()

# New column header and value for row 2 (the single data row).
$ws.Range("F1").Value = "Previous balance"
$ws.Range("F2").Value = 200000

# The "previous balance" figure gets its own (distinct) cell style/font.
$ws.Range("F2").Font.Name = "Arial"
$ws.Range("F2").Font.Size = 10
$ws.Range("F2").Font.Color = 0

# Corrections to the balance-drop related flags, shifted one column right
# after the insert above.
$ws.Range("M2").Value = "yes"   # internet banking
$ws.Range("P2").Value = "yes"   # neft_rtgs

$ws.Range("F4").Select()
